$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every data row
# (rows 2-118). The automatic update bumps this date by one day,
# from serial 45181 (2023-09-12) to serial 45182 (2023-09-13) for all rows.
$ws.Range("C2:C118").Value = 45182
